$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.74580449258613424
$ws.Range("D2").Value = 0.70253011264985621
$ws.Range("A3").Value = 0.69900658042331953
$ws.Range("D3").Value = 0.85993051724872793
$ws.Range("AP3").Value = 0.94160148344764671
$ws.Range("I4").Value = 0.85676225756284119
$ws.Range("AH4").Value = 0.79035482555033032
$ws.Range("C5").Value = 0.80692410207146459
$ws.Range("F5").Value = 0.96678221415013377
$ws.Range("G5").Value = 0.69831386309041976
$ws.Range("H6").Value = 0.8269879526374857
$ws.Range("F7").Value = 0.82212990654616136
$ws.Range("I7").Value = 0.94895500074408068
$ws.Range("X7").Value = 0.7367315684821889
$ws.Range("I8").Value = 0.98560102744511424
$ws.Range("I10").Value = 0.83182420308100491
$ws.Range("AJ10").Value = 0.86459042171974709
$ws.Range("I11").Value = 0.89957469875771678
$ws.Range("T11").Value = 0.86091497941674189
$ws.Range("J12").Value = 0.76447448743919399
$ws.Range("K13").Value = 0.86464417518542525
$ws.Range("L13").Value = 0.77239351944830181
$ws.Range("N13").Value = 0.92619090211947497
$ws.Range("L14").Value = 0.56198928441095575
$ws.Range("U14").Value = 0.77169484713688985
$ws.Range("BD14").Value = 0.73749851679986511
$ws.Range("M15").Value = 0.81557831153191285
$ws.Range("N15").Value = 0.93131009464369097
$ws.Range("Q15").Value = 0.70644983821748175
$ws.Range("AI15").Value = 0.70556157792259855
$ws.Range("O16").Value = 0.88992424764877442
$ws.Range("Q16").Value = 0.99801966628180383
$ws.Range("P18").Value = 0.68072205171695321
$ws.Range("Q18").Value = 0.95201936152276967
$ws.Range("S18").Value = 0.8234156956718508
$ws.Range("AC18").Value = 0.76294619849066114
$ws.Range("AY18").Value = 0.96679066325086982
$ws.Range("AQ19").Value = 0.84917070133750316
$ws.Range("S20").Value = 0.94660570868544203
$ws.Range("W21").Value = 0.73558221251151745
$ws.Range("AR21").Value = 0.65349859452801151
$ws.Range("T22").Value = 0.79745181144119148
$ws.Range("U22").Value = 0.99946276433356496
$ws.Range("BO22").Value = 0.75134467784862657
$ws.Range("V23").Value = 0.69231445566426819
$ws.Range("Y23").Value = 0.94226201190492676
$ws.Range("W24").Value = 0.9509525811319548
$ws.Range("AA25").Value = 0.84740881529010115
$ws.Range("U26").Value = 0.97722884187488501
$ws.Range("X26").Value = 0.88179617940622168
$ws.Range("AF26").Value = 0.89072404469202104
$ws.Range("AC27").Value = 0.88307064898674437
$ws.Range("BI27").Value = 0.72376797185369024
$ws.Range("Z28").Value = 0.96749052331073249
$ws.Range("AB29").Value = 0.77693749502622667
$ws.Range("AB30").Value = 0.74503117936666341
$ws.Range("AE30").Value = 0.9002842472340804
$ws.Range("AF30").Value = 0.92826623201465197
$ws.Range("AC31").Value = 0.93039443666610522
$ws.Range("Q32").Value = 0.92532072065235549
$ws.Range("R32").Value = 0.94216271665299134
$ws.Range("AE32").Value = 0.80737905749142114
$ws.Range("AE33").Value = 0.80271394130881824
$ws.Range("AH33").Value = 0.96503880880591431
$ws.Range("AI33").Value = 0.80778073473917034
$ws.Range("AF34").Value = 0.69067566387027224
$ws.Range("BK34").Value = 0.66962040177795179
$ws.Range("T35").Value = 0.77429649378833931
$ws.Range("AH35").Value = 0.69138491012751402
$ws.Range("AJ35").Value = 0.88453541805897706
$ws.Range("AK36").Value = 0.98248710830698105
$ws.Range("AL36").Value = 0.64839831192108743
$ws.Range("AM37").Value = 0.87410584502245947
$ws.Range("AK38").Value = 0.83323038290837159
$ws.Range("AN38").Value = 0.67168467186470826
$ws.Range("AL39").Value = 0.96978740466886848
$ws.Range("AN39").Value = 0.52738198970537919
$ws.Range("AO39").Value = 0.88482221720721643
$ws.Range("D40").Value = 0.60923536564821656
$ws.Range("AP40").Value = 0.8788377921987458
$ws.Range("AQ41").Value = 0.68424226224261675
$ws.Range("G42").Value = 0.80456006316647755
$ws.Range("AO42").Value = 0.72382769921823242
$ws.Range("AR42").Value = 0.69072745748646458
$ws.Range("AR43").Value = 0.66017184202215384
$ws.Range("AT44").Value = 0.87881525475073174
$ws.Range("AQ45").Value = 0.87412001519178584
$ws.Range("AS46").Value = 0.96452943608013086
$ws.Range("AU46").Value = 0.7985352732366886
$ws.Range("AS47").Value = 0.910442205782799
$ws.Range("AW47").Value = 0.8329568998997765
$ws.Range("AU48").Value = 0.79137345790320457
$ws.Range("AW48").Value = 0.71884579862766862
$ws.Range("AX48").Value = 0.93054179848031227
$ws.Range("A49").Value = 0.86312688452600028
$ws.Range("S49").Value = 0.69719401634686062
$ws.Range("AZ50").Value = 0.72670492296802669
$ws.Range("AA51").Value = 0.7908867667172409
$ws.Range("AX51").Value = 0.66203949054198619
$ws.Range("AV52").Value = 0.8830391243510014
$ws.Range("BA52").Value = 0.91033423026424876
$ws.Range("C53").Value = 0.97991875433837206
$ws.Range("AY53").Value = 0.89127488317498016
$ws.Range("BB53").Value = 0.92880949435663041
$ws.Range("AT54").Value = 0.96577354822060457
$ws.Range("BD54").Value = 0.70681323659601136
$ws.Range("N55").Value = 0.97550015736651297
$ws.Range("BB55").Value = 0.85239632698488443
$ws.Range("BE55").Value = 0.70557880664135197
$ws.Range("Y56").Value = 0.92031017237837554
$ws.Range("AS57").Value = 0.71093578297428239
$ws.Range("BG57").Value = 0.98948926315487196
$ws.Range("G58").Value = 0.92868898056175231
$ws.Range("BD58").Value = 0.86295943088825466
$ws.Range("BG58").Value = 0.82912970019497578
$ws.Range("A60").Value = 0.94551878715207327
$ws.Range("BG60").Value = 0.78926501768814483
$ws.Range("H61").Value = 0.89448316387785387
$ws.Range("BG61").Value = 0.84248574422234679
$ws.Range("BH61").Value = 0.91833976271714035
$ws.Range("BJ61").Value = 0.85336805693373297
$ws.Range("BA62").Value = 0.74140213140618283
$ws.Range("BH62").Value = 0.92414245933931638
$ws.Range("BK62").Value = 0.88471560536540839
$ws.Range("BP62").Value = 0.67103517279574021
$ws.Range("BM63").Value = 0.95103415902647637
$ws.Range("I64").Value = 0.7362329976898272
$ws.Range("AK64").Value = 0.8806120567209893
$ws.Range("BK64").Value = 0.63110523106824024
$ws.Range("BM66").Value = 0.78648572393827332
$ws.Range("A67").Value = 0.59297949411845541
$ws.Range("BM67").Value = 0.69402098742427887
$ws.Range("BN67").Value = 0.92853555076715644
$ws.Range("BP67").Value = 0.91582422014789566
$ws.Range("BF68").Value = 0.77311261371490203
$ws.Range("BN68").Value = 0.75489731425047046
